# Rename the img* sheets: old name X-img -> img-X (e.g. himg -> imgh)
$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("himg").Name = "imgh"
$wb.Worksheets.Item("timg").Name = "imgt"
$wb.Worksheets.Item("simg").Name = "imgs"
$wb.Worksheets.Item("gimg").Name = "imgg"
$wb.Worksheets.Item("wimg").Name = "imgw"
$wb.Worksheets.Item("bimg").Name = "imgb"
$wb.Worksheets.Item("eimg").Name = "imge"

# Restore / set the active cell selection on every sheet (single-cell
# selections, dropping the old secondary "C16:D17" range reference).
$wb.Worksheets.Item("expert").Select()
$wb.Worksheets.Item("expert").Range("A3").Select()

$wb.Worksheets.Item("task").Select()
$wb.Worksheets.Item("task").Range("B4").Select()

$wb.Worksheets.Item("assign").Select()
$wb.Worksheets.Item("assign").Range("C16").Select()

$wb.Worksheets.Item("xbday").Select()
$wb.Worksheets.Item("xbday").Range("E5").Select()

$wb.Worksheets.Item("ubday").Select()
$wb.Worksheets.Item("ubday").Range("C2").Select()

$wb.Worksheets.Item("ebday").Select()
$wb.Worksheets.Item("ebday").Range("G5").Select()

$wb.Worksheets.Item("period").Select()
$wb.Worksheets.Item("period").Range("C2").Select()

$wb.Worksheets.Item("pbsum").Select()
$wb.Worksheets.Item("pbsum").Range("H16").Select()

$wb.Worksheets.Item("holiday").Select()
$wb.Worksheets.Item("holiday").Range("F13").Select()

$wb.Worksheets.Item("misc").Select()
$wb.Worksheets.Item("misc").Range("H10").Select()

$wb.Worksheets.Item("imgh").Select()
$wb.Worksheets.Item("imgh").Range("F3").Select()

$wb.Worksheets.Item("imgt").Select()
$wb.Worksheets.Item("imgt").Range("F3").Select()

$wb.Worksheets.Item("imgs").Select()
$wb.Worksheets.Item("imgs").Range("F3").Select()

$wb.Worksheets.Item("imgg").Select()
$wb.Worksheets.Item("imgg").Range("H1").Select()

$wb.Worksheets.Item("imgw").Select()
$wb.Worksheets.Item("imgw").Range("A1").Select()

$wb.Worksheets.Item("imgb").Select()
$wb.Worksheets.Item("imgb").Range("A1").Select()

# "imge" (formerly "eimg") is the final active sheet/tab in the target
# workbook (activeTab=16, tabSelected on this sheet).
$wb.Worksheets.Item("imge").Select()
$wb.Worksheets.Item("imge").Range("A1").Select()
